# Restore C10 ("R30" rule, "Integer min" / From column) from 18 to 1,
# matching revision f135250c7393bc5c76059d56caf7e133ff65dd5e.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 1
